# "got dim 0 working with adjacency lists and k_n"
#
# The workbook tracks MST statistics for several dimensions (2, 3, 4 already
# present on rows 10-26). This adds the "dimension 0" block (computed with
# the new adjacency-list implementation) as rows 28-32, extends a couple of
# the existing dimension-0 sample rows (5-7) and the row-23/24 n=2048 column
# with new adjacency-list figures, and extends the scratch table at the
# bottom of the sheet (rows 33/34/38-40) with the same new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- extend the existing dimension-0 quick table (rows 5-7) with the new
#     n = 2048 .. 16834 columns (L:O) ---------------------------------------
$ws.Range("L5").Value = 1.198184
$ws.Range("M5").Value = 1.207412
$ws.Range("N5").Value = 1.196903
$ws.Range("O5").Value = 1.207155

$ws.Range("L6").Value = 0.00489
$ws.Range("M6").Value = 0.002424
$ws.Range("N6").Value = 0.001778
$ws.Range("O6").Value = 0.000677

$ws.Range("O7").Value = 5

# --- dimension 4 block (rows 22-26): add the n = 2048 trial column --------
$ws.Range("L23").Value = 1642.531494
$ws.Range("L24").Value = 1.720878

# --- new "dimension 0, adjacency list" block: rows 28-32 -------------------
$ws.Range("A28").Value = "NEW WITH ADJACENCY LIST"
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = "n"
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 32
$ws.Range("G28").Value = 64
$ws.Range("H28").Value = 128
$ws.Range("I28").Value = 256
$ws.Range("J28").Value = 512
$ws.Range("K28").Value = 1024
$ws.Range("L28").Value = 2048
$ws.Range("M28").Value = 4096
$ws.Range("N28").Value = 8192
$ws.Range("O28").Value = 16834
$ws.Range("P28").Value = 32768
$ws.Range("Q28").Value = 65536

$ws.Range("D29").Value = "Average MST Weight"
$ws.Range("E29").Value = 1.318274
$ws.Range("F29").Value = 1.481506
$ws.Range("G29").Value = 1.488744
$ws.Range("H29").Value = 1.777309
$ws.Range("I29").Value = 1.792974
$ws.Range("J29").Value = 2.002005
$ws.Range("K29").Value = 2.104817
$ws.Range("L29").Value = 1.198184
$ws.Range("M29").Value = 1.207412
$ws.Range("N29").Value = 1.196903
$ws.Range("O29").Value = 1.207155

$ws.Range("D30").Value = "Max included edge"
$ws.Range("E30").Value = 0.543776
$ws.Range("F30").Value = 0.36075
$ws.Range("G30").Value = 0.210741
$ws.Range("H30").Value = 0.098394
$ws.Range("I30").Value = 0.051644
$ws.Range("J30").Value = 0.02648
$ws.Range("K30").Value = 0.014001
$ws.Range("L30").Value = 0.00489
$ws.Range("M30").Value = 0.002424
$ws.Range("N30").Value = 0.001778
$ws.Range("O30").Value = 0.000677

$ws.Range("D31").Value = "(trials)"
$ws.Range("E31").Value = 10000
$ws.Range("F31").Value = 10000
$ws.Range("G31").Value = 10000
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("O31").Value = 5

$ws.Range("C32").Value = "k(n)"
$ws.Range("D32").Value = "0.025+(1.25)*4.53261/(n^0.754872)"
$ws.Range("E32").Formula = "=0.025+(1.25)*4.53261/POWER(E28,0.754872)"
# F32:K32 share one formula (relative refs adjust per column), same as the
# existing F8:K8 block above.
$ws.Range("F32:K32").Formula = "=0.025+(1.25)*4.53261/POWER(F28,0.754872)"

# --- scratch table at bottom of sheet: extend with the n=2048..16834 cols -
$ws.Range("O38").Value = 1.198184
$ws.Range("P38").Value = 1.207412
$ws.Range("Q38").Value = 1.196903
$ws.Range("R38").Value = 1.207155

$ws.Range("O39").Value = 0.00489
$ws.Range("P39").Value = 0.002424
$ws.Range("Q39").Value = 0.001778
$ws.Range("R39").Value = 0.000677

# --- move the active selection to where the new work was done -------------
$ws.Range("H18").Select()
